# Regenerate the "K" column (strikeouts) values for each start in the
# save_data sheet, replacing the previous Strike# derived figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 2
    10 = 2
    11 = 5
    12 = 3
    13 = 2
    14 = 5
    15 = 2
    16 = 3
    17 = 2
    18 = 2
    19 = 2
    20 = 0
    21 = 1
    22 = 1
    23 = 2
    24 = 0
    25 = 2
    26 = 2
    27 = 1
    28 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
